$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume data (GitHub Actions scrape refresh).
# Price column (D) values are plain text in the source data (e.g. "1.002",
# "0.5120", "25.900.26"), so force text format before assigning to avoid Excel
# auto-converting dotted numeric-looking strings into actual numbers and
# dropping significant trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.900.26"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.632.60"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.05"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5120"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2573"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06343"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.46"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07782"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.281"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.636.56"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.858.84"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5508"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.91"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₅7645"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.942.36"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.421"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.72"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.864"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.037"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.889"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.79"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1256"
$ws.Range("E27").Value = "  +5.05%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.755"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.55"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04876"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.241"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.188"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.541"
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.374"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8978"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5518"
$ws.Range("E37").Value = "  +1.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.541"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.117.44"
$ws.Range("E39").Value = "  -2.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01558"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.587"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7971"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.766.77"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈118"
$ws.Range("E46").Value = "  -6.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4448"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.74"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05131"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.591"
$ws.Range("E51").Value = "  +3.21%  "
